$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: update Fecha, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg
$ws.Range("D2").Value = 44533
$ws.Range("N2").Value = 16000
$ws.Range("O2").Value = 17000
$ws.Range("P2").Value = 16500
$ws.Range("S2").Value = 825

# Row 4: update Fecha, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg
$ws.Range("D4").Value = 44357
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 14500
$ws.Range("S4").Value = 725
